# Update row 8 (year 2025) metrics in the recorrencia anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1406
$ws.Range("D8").Value = 214
$ws.Range("E8").Value = 1192
$ws.Range("F8").Value = 8.777686628383922
$ws.Range("G8").Value = 84.77951635846372
$ws.Range("H8").Value = 15.22048364153627
